$wb = $excel.ActiveWorkbook

# --- Sheet2: move the small "any heading"/"Amount" lookup table from J3:K4
#     down to F11:G12, directly under the first mini table (same heading
#     column as the other table on the sheet) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("J3:K4").Cut($ws2.Range("F11"))
$ws2.Range("J3:K4").Clear()
$ws2.Columns.Item(6).AutoFit()
$ws2.Range("J6").Select()
$ws2.PageSetup.Orientation = 1

# --- Sheet3: same change ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("J3:K4").Cut($ws3.Range("F11"))
$ws3.Range("J3:K4").Clear()
$ws3.Columns.Item(6).AutoFit()
$ws3.Range("G15").Select()
